$d = $word.ActiveDocument

# Curly right single quotation mark (U+2019) used in "man's"
$apos = [char]0x2019
$newText = " The man" + $apos + "s goal is to get everything to the other side of the river"

# Locate the paragraph that ends with "...one more." (the "Cat, Parrot, and Bag
# of Seed" answer paragraph) by searching for the distinctive trailing phrase.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "one more\.\s*$") {
        $target = $p
        break
    }
}

$r = $target.Range

# The paragraph Range.End sits one past the paragraph mark; back up one
# character so the new text is inserted right after the final "."  and
# before the paragraph mark (and ahead of the _GoBack bookmark).
$insPos = $r.End - 1
$insPoint = $d.Range($insPos, $insPos)
$insPoint.InsertAfter($newText)

# Re-select just the newly inserted text so we can stamp it with explicit
# run formatting. Toggling the size briefly forces Word to start a fresh
# run (rather than silently merging back into the preceding identical-
# format run) so the inserted text keeps its own <w:r> in the OOXML,
# matching the source edit.
$newRange = $d.Range($insPos, $insPos + $newText.Length)
$newRange.Font.Name = "Times New Roman"
$newRange.Font.Size = 10
$newRange.Font.Name = "Times New Roman"
$newRange.Font.Size = 12
